$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 660
$ws.Range("I49").Value = 430
$ws.Range("J49").Value = 775
$ws.Range("K49").Value = 1290
$ws.Range("L49").Value = 2325
$ws.Range("M49").Value = -1154
$ws.Range("N49").Value = -2597
$ws.Range("H51").Value = 5855.7144
$ws.Range("J51").Value = 3498
$ws.Range("L51").Value = 3498
$ws.Range("N51").Value = -4466
$ws.Range("H62").Value = 6026.316
$ws.Range("I62").Value = 5163.4546
$ws.Range("K62").Value = 5163.4546
$ws.Range("M62").Value = -4539.4546
$ws.Range("H65").Value = 6026.316
$ws.Range("I65").Value = 5163.4546
$ws.Range("K65").Value = 25817.273
$ws.Range("M65").Value = -22697.273
$ws.Range("H86").Value = 12893.889
$ws.Range("I86").Value = 1771.8572
$ws.Range("K86").Value = 1771.8572
$ws.Range("M86").Value = -648.8571999999999
$ws.Range("H89").Value = 12893.889
$ws.Range("I89").Value = 1771.8572
$ws.Range("K89").Value = 8859.286
$ws.Range("M89").Value = -3243.286
$ws.Range("H98").Value = 688.82355
$ws.Range("I98").Value = 644.375
$ws.Range("J98").Value = 1400
$ws.Range("K98").Value = 644.375
$ws.Range("L98").Value = 1400
$ws.Range("M98").Value = 853.625
$ws.Range("N98").Value = -4396
$ws.Range("H113").Value = 58827092
$ws.Range("I113").Value = 142858380
$ws.Range("J113").Value = 5180.7
$ws.Range("K113").Value = 142858380
$ws.Range("L113").Value = 5180.7
$ws.Range("M113").Value = -142855126
$ws.Range("N113").Value = -11688.7
$ws.Range("H116").Value = 11765420
$ws.Range("I116").Value = 40324524
$ws.Range("J116").Value = 5788.5884
$ws.Range("K116").Value = 40324524
$ws.Range("L116").Value = 5788.5884
$ws.Range("M116").Value = -40321082
$ws.Range("N116").Value = -12672.5884
$ws.Range("H122").Value = 688.82355
$ws.Range("I122").Value = 644.375
$ws.Range("J122").Value = 1400
$ws.Range("K122").Value = 1933.125
$ws.Range("L122").Value = 4200
$ws.Range("M122").Value = 516.875
$ws.Range("N122").Value = -9100
$ws.Range("H127").Value = 1160.5385
$ws.Range("J127").Value = 1820
$ws.Range("L127").Value = 5460
$ws.Range("N127").Value = -15380
$ws.Range("H129").Value = 271179.22
$ws.Range("I129").Value = 244.25
$ws.Range("J129").Value = 304019.8
$ws.Range("K129").Value = 732.75
$ws.Range("L129").Value = 912059.3999999999
$ws.Range("M129").Value = 4267.25
$ws.Range("N129").Value = -922059.3999999999
$ws.Range("H137").Value = 87416.66
$ws.Range("I137").Value = 107097.37
$ws.Range("J137").Value = 4320.3335
$ws.Range("K137").Value = 321292.11
$ws.Range("L137").Value = 12961.0005
$ws.Range("M137").Value = -318742.11
$ws.Range("N137").Value = -18061.0005
$ws.Range("H138").Value = 2976.4314
$ws.Range("I138").Value = 2142.85
$ws.Range("J138").Value = 3514.2258
$ws.Range("K138").Value = 6428.549999999999
$ws.Range("L138").Value = 10542.6774
$ws.Range("M138").Value = -1288.549999999999
$ws.Range("N138").Value = -20822.6774

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7576.663
$ws.Range("I32").Value = 5407.393
$ws.Range("J32").Value = 24142
$ws.Range("K32").Value = 5407.393
$ws.Range("L32").Value = 24142
$ws.Range("M32").Value = -5120.393
$ws.Range("N32").Value = -24716
$ws.Range("H45").Value = 2663.3225
$ws.Range("I45").Value = 2193.4
$ws.Range("K45").Value = 2193.4
$ws.Range("M45").Value = -1816.4
$ws.Range("H74").Value = 25001508
$ws.Range("I74").Value = 34483410
$ws.Range("K74").Value = 34483410
$ws.Range("M74").Value = -34482536
$ws.Range("H77").Value = 25001508
$ws.Range("I77").Value = 34483410
$ws.Range("K77").Value = 172417050
$ws.Range("M77").Value = -172412682
$ws.Range("H110").Value = 1455.5333
$ws.Range("I110").Value = 1009.7273
$ws.Range("J110").Value = 2681.5
$ws.Range("K110").Value = 1009.7273
$ws.Range("L110").Value = 2681.5
$ws.Range("M110").Value = 1035.2727
$ws.Range("N110").Value = -6771.5
$ws.Range("H122").Value = 3769.111
$ws.Range("I122").Value = 3802.875
$ws.Range("K122").Value = 11408.625
$ws.Range("M122").Value = -8958.625
$ws.Range("H132").Value = 8074793
$ws.Range("I132").Value = 10205870
$ws.Range("J132").Value = 42272.31
$ws.Range("K132").Value = 30617610
$ws.Range("L132").Value = 126816.93
$ws.Range("M132").Value = -30615080
$ws.Range("N132").Value = -131876.93

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 7657.2856
$ws.Range("I134").Value = 8679.643
$ws.Range("J134").Value = 5612.5713
$ws.Range("K134").Value = 26038.929
$ws.Range("L134").Value = 16837.7139
$ws.Range("M134").Value = -23503.929
$ws.Range("N134").Value = -21907.7139

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 30261.285
$ws.Range("J52").Value = 30261.285
$ws.Range("L52").Value = 30261.285
$ws.Range("N52").Value = -30849.285
$ws.Range("H141").Value = 25947.334
$ws.Range("J141").Value = 25947.334
$ws.Range("L141").Value = 25947.334
$ws.Range("N141").Value = -36307.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 244.15
$ws.Range("I44").Value = 236.26666
$ws.Range("J44").Value = 267.8
$ws.Range("K44").Value = 708.79998
$ws.Range("L44").Value = 803.4000000000001
$ws.Range("M44").Value = -310.79998
$ws.Range("N44").Value = -1599.4
$ws.Range("H122").Value = 1429.96
$ws.Range("J122").Value = 1429.96
$ws.Range("L122").Value = 12869.64
$ws.Range("N122").Value = -17769.64
$ws.Range("H131").Value = 688.11
$ws.Range("J131").Value = 723.0449
$ws.Range("L131").Value = 2169.1347
$ws.Range("N131").Value = -12249.1347
$ws.Range("H139").Value = 3553.16
$ws.Range("I139").Value = 2313.2856
$ws.Range("J139").Value = 5131.1816
$ws.Range("K139").Value = 6939.8568
$ws.Range("L139").Value = 15393.5448
$ws.Range("M139").Value = -1799.8568
$ws.Range("N139").Value = -25673.5448

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 25183.334
$ws.Range("J57").Value = 29920
$ws.Range("L57").Value = 29920
$ws.Range("N57").Value = -31560
$ws.Range("H102").Value = 5404.1113
$ws.Range("I102").Value = 4827.875
$ws.Range("K102").Value = 4827.875
$ws.Range("M102").Value = -3205.875
$ws.Range("H126").Value = 3918.182
$ws.Range("I126").Value = 2785.7144
$ws.Range("K126").Value = 8357.143199999999
$ws.Range("M126").Value = -5887.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 2166.6667
$ws.Range("J24").Value = 2166.6667
$ws.Range("L24").Value = 2166.6667
$ws.Range("N24").Value = -2852.6667
$ws.Range("H25").Value = 7000
$ws.Range("J25").Value = 7000
$ws.Range("L25").Value = 7000
$ws.Range("N25").Value = -7460
$ws.Range("H40").Value = 56755.816
$ws.Range("I40").Value = 80708.664
$ws.Range("J40").Value = 5428.2856
$ws.Range("K40").Value = 80708.664
$ws.Range("L40").Value = 5428.2856
$ws.Range("M40").Value = -80572.664
$ws.Range("N40").Value = -5700.2856
$ws.Range("H46").Value = 1349.75
$ws.Range("I46").Value = 999
$ws.Range("J46").Value = 1466.6666
$ws.Range("K46").Value = 999
$ws.Range("L46").Value = 1466.6666
$ws.Range("M46").Value = -811
$ws.Range("N46").Value = -1842.6666
$ws.Range("H93").Value = 2031.25
$ws.Range("I93").Value = 2115.3845
$ws.Range("J93").Value = 1666.6666
$ws.Range("K93").Value = 2115.3845
$ws.Range("L93").Value = 1666.6666
$ws.Range("M93").Value = -867.3845000000001
$ws.Range("N93").Value = -4162.6666
$ws.Range("H122").Value = 1311106.2
$ws.Range("I122").Value = 1637341.1
$ws.Range("J122").Value = 6166.6665
$ws.Range("K122").Value = 4912023.300000001
$ws.Range("L122").Value = 18499.9995
$ws.Range("M122").Value = -4909573.300000001
$ws.Range("N122").Value = -23399.9995
$ws.Range("H132").Value = 3005.3076
$ws.Range("I132").Value = 2397.3044
$ws.Range("K132").Value = 7191.9132
$ws.Range("M132").Value = -4661.9132
$ws.Range("H136").Value = 2922.3333
$ws.Range("I136").Value = 2969
$ws.Range("K136").Value = 8907
$ws.Range("M136").Value = -6357

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 5042.2
$ws.Range("I20").Value = 3000
$ws.Range("K20").Value = 3000
$ws.Range("M20").Value = -2760
$ws.Range("H62").Value = 4640.2856
$ws.Range("I62").Value = 3501
$ws.Range("J62").Value = 5096
$ws.Range("K62").Value = 3501
$ws.Range("L62").Value = 5096
$ws.Range("M62").Value = -2877
$ws.Range("N62").Value = -6344
$ws.Range("H65").Value = 4640.2856
$ws.Range("I65").Value = 3501
$ws.Range("J65").Value = 5096
$ws.Range("K65").Value = 17505
$ws.Range("L65").Value = 25480
$ws.Range("M65").Value = -14385
$ws.Range("N65").Value = -31720
$ws.Range("H122").Value = 1568.3125
$ws.Range("J122").Value = 1311.625
$ws.Range("L122").Value = 3934.875
$ws.Range("N122").Value = -8834.875
$ws.Range("H126").Value = 1998.8636
$ws.Range("I126").Value = 1529.5
$ws.Range("K126").Value = 4588.5
$ws.Range("M126").Value = -2118.5
$ws.Range("H132").Value = 11905944
$ws.Range("I132").Value = 15625825
$ws.Range("K132").Value = 46877475
$ws.Range("M132").Value = -46874945
